$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header "Comments" -> "Notes" in O1
$ws.Range("O1").Value = "Notes"

# Update the table column name for the "Comments" column
$table = $ws.ListObjects.Item(1)
$table.ListColumns.Item("Comments").Name = "Notes"

# Update selection / scroll position
$ws.Range("O2").Select()
$excel.ActiveWindow.ScrollColumn = 7
